$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3,0,3,3),
    @(5,2,4,0),
    @(3,1,4,2),
    @(5,1,5,2),
    @(5,0,7,2),
    @(6,0,6,2),
    @(3,2,5,0),
    @(4,3,3,0),
    @(5,2,4,0),
    @(5,3,4,0),
    @(5,0,5,2),
    @(4,3,4,0),
    @(5,2,5,1),
    @(4,1,4,2),
    @(5,3,4,0),
    @(3,2,3,1),
    @(3,2,4,1),
    @(6,0,5,2),
    @(7,2,6,0),
    @(3,0,4,3),
    @(3,2,5,1),
    @(4,0,4,2),
    @(5,2,5,0),
    @(3,3,3,0),
    @(5,3,3,0)
)

$startRow = 2176
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$endRow = $startRow + $data.Count - 1

# Update the selection to mirror what Excel leaves after pasting the new
# rows in (one past the last data row), and nudge the viewport so the
# freshly entered rows are in view.
$ws.Range("A" + ($endRow + 1)).Select()
try {
    $excel.ActiveWindow.ScrollRow = 2173
} catch {
    # ActiveWindow scroll position isn't always addressable; ignore.
}
